$wb = $excel.ActiveWorkbook

# --- Add "Problem 8" worksheet (product/price table) after the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws8 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws8.Name = "Problem 8"

$ws8.Range("A1").Value = "Product_name"
$ws8.Range("B1").Value = "Price"

$ws8.Range("A2").Value = "Apple"
$ws8.Range("A3").Value = "Banana"
$ws8.Range("A4").Value = "Orange"
$ws8.Range("A5").Value = "Mango"

$ws8.Range("B2").Value = 50
$ws8.Range("B3").Value = 30
$ws8.Range("B4").Value = 70
$ws8.Range("B5").Value = 40

$ws8.Columns.Item(1).AutoFit()
$ws8.Range("B6").Select()

# --- Add "Problem 9" worksheet (employee table) after "Problem 8" ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws9 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet2)
$ws9.Name = "Problem 9"

$ws9.Range("A1").Value = "emp_id"
$ws9.Range("B1").Value = "first"
$ws9.Range("C1").Value = "last"
$ws9.Range("D1").Value = "years_exp"

$ws9.Range("A2").Value = 101
$ws9.Range("A3").Value = 102
$ws9.Range("A4").Value = 103
$ws9.Range("A5").Value = 104

$ws9.Range("B2").Value = "John"
$ws9.Range("B3").Value = "Jane"
$ws9.Range("B4").Value = "Emily"
$ws9.Range("B5").Value = "Michael"

$ws9.Range("C2").Value = "Doe"
$ws9.Range("C3").Value = "Smith"
$ws9.Range("C4").Value = "Johnson"
$ws9.Range("C5").Value = "Williams"

$ws9.Range("D2").Value = 5
$ws9.Range("D3").Value = 10
$ws9.Range("D4").Value = 3
$ws9.Range("D5").Value = 8

$ws9.Range("D24").Select()
